# NIT-9008151452.xlsx — "Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta"
#
# The worker LUIS FERNANDO MIRANDA GUTIERREZ (CC 1143127438) previously had a
# single overdue period (2507) listed on the "Estado de Cuenta" sheet. This
# edit adds a second overdue period (2508) for the same worker as a new data
# row, then updates the summary totals ("VALOR MORA" and "Cant. Periodos")
# to reflect the extra period.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new (blank) row right below the worker's existing data row (16).
# This pushes the signature block (old rows 21-22) down to rows 22-23,
# matching the target layout where the new row becomes row 17.
$ws.Rows.Item(17).Insert()

# Give the new row the same look as row 16 (borders, fonts, number formats)
# by copying row 16's formatting onto it.
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J17").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Populate the new row with the same worker, but for the next overdue
# period: 2507 -> 2508. Everything else (doc type, doc number, name, value
# mora, salario basico) is identical to row 16.
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1143127438"
$ws.Range("D17").Value = "LUIS FERNANDO MIRANDA GUTIERREZ"
$ws.Range("E17").Value = "2508"
$ws.Range("F17").Value = 166480
$ws.Range("G17").Value = 4162000

# Update the header totals: two overdue periods now exist for this worker,
# so the overdue amount doubles (166480 -> 332960) and the period count
# goes from 1 to 2.
$ws.Range("E11").Value = 332960
$ws.Range("F13").Value = 2
